$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.038.99"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").Value = "3.418.16"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("D5").Value = "406.15"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").Value = "131.52"
$ws.Range("E6").Value = "  +1.59%  "
$ws.Range("E7").Value = "  -2.04%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "0.688"
$ws.Range("E9").Value = "  +1.95%  "
$ws.Range("D10").Value = "0.134"
$ws.Range("E10").Value = "  +4.93%  "
$ws.Range("D11").Value = "41.75"
$ws.Range("E11").Value = "  -1.56%  "
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "19.84"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").Value = "8.38"
$ws.Range("E14").Value = "  -2.32%  "
$ws.Range("D15").Value = "3.379.56"
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "61.926.54"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "11.56"
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").Value = "0.0000147"
$ws.Range("E19").Value = "  +9.27%  "
$ws.Range("D20").Value = "3.16"
$ws.Range("E20").Value = "  -2.64%  "
$ws.Range("D21").Value = "83.84"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").Value = "313.54"
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("D23").Value = "12.75"
$ws.Range("E23").Value = "  -2.68%  "
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").Value = "4.76"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").Value = "29.61"
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("D27").Value = "'7.90"
$ws.Range("E27").Value = "  +5.16%  "
$ws.Range("D28").Value = "8.12"
$ws.Range("E28").Value = "  -5.95%  "
$ws.Range("D29").Value = "2.73"
$ws.Range("E29").Value = "  +5.20%  "
$ws.Range("D30").Value = "0.173"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").Value = "43.48"
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("D32").Value = "0.115"
$ws.Range("E32").Value = "  -1.20%  "
$ws.Range("D33").Value = "'11.30"
$ws.Range("E33").Value = "  -3.85%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("D35").Value = "0.0486"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("D36").Value = "51.55"
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("D39").Value = "3.34"
$ws.Range("E39").Value = "  -2.84%  "
$ws.Range("D40").Value = "0.313"
$ws.Range("E40").Value = "  +9.51%  "
$ws.Range("D41").Value = "139.62"
$ws.Range("E41").Value = "  +2.69%  "
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").Value = "1.97"
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("D44").Value = "3.92"
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("D45").Value = "16.71"
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").Value = "21.23"
$ws.Range("E47").Value = "  -3.39%  "
$ws.Range("D48").Value = "2.100.87"
$ws.Range("E48").Value = "  -2.44%  "
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").Value = "1.92"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "1.69"
$ws.Range("E51").Value = "  +14.72%  "
